# Fill in the room ("salle") column (F) for each course occurrence row.
# "FSQTEL - C" (cours) sessions -> "U3-Amphi"; the TD/controle sessions -> "U3-110".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "U3-Amphi"
$ws.Range("F5").Value = "U3-110"
$ws.Range("F7").Value = "U3-110"
$ws.Range("F10").Value = "U3-110"
$ws.Range("F11").Value = "U3-110"
$ws.Range("F14").Value = "U3-110"
$ws.Range("F17").Value = "U3-Amphi"
